$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (preserve string type / leading/trailing formatting like "305.30" or "14")
# by temporarily applying a text number format to the affected range before writing,
# then restore the original "Normal" style so no formatting diff is introduced.
$editRange = $ws.Range("D2:G51")
$editRange.NumberFormat = "@"

$ws.Range("D2").Value = "305.30"
$ws.Range("E2").Value = "0.26%"
$ws.Range("G2").Value = "14"
$ws.Range("D3").Value = "38.26"
$ws.Range("E3").Value = "6.99%"
$ws.Range("G3").Value = "14"
$ws.Range("D4").Value = "5.094"
$ws.Range("E4").Value = "1.20%"
$ws.Range("G4").Value = "14"
$ws.Range("D5").Value = "0.08056"
$ws.Range("E5").Value = "1.07%"
$ws.Range("G5").Value = "14"
$ws.Range("D6").Value = "1.936"
$ws.Range("E6").Value = "4.05%"
$ws.Range("G6").Value = "14"
$ws.Range("D7").Value = "4.195"
$ws.Range("E7").Value = "1.75%"
$ws.Range("G7").Value = "14"
$ws.Range("D8").Value = "7.946"
$ws.Range("E8").Value = "2.20%"
$ws.Range("G8").Value = "14"
$ws.Range("D9").Value = "0.9281"
$ws.Range("E9").Value = "0.75%"
$ws.Range("G9").Value = "14"
$ws.Range("D10").Value = "0.1434"
$ws.Range("E10").Value = "12.59%"
$ws.Range("G10").Value = "14"
$ws.Range("D11").Value = "0.1921"
$ws.Range("E11").Value = "1.88%"
$ws.Range("G11").Value = "14"
$ws.Range("D12").Value = "0.09021"
$ws.Range("E12").Value = "-0.39%"
$ws.Range("G12").Value = "14"
$ws.Range("D13").Value = "0.03496"
$ws.Range("E13").Value = "2.23%"
$ws.Range("G13").Value = "14"
$ws.Range("D14").Value = "0.09774"
$ws.Range("E14").Value = "-0.78%"
$ws.Range("G14").Value = "14"
$ws.Range("D15").Value = "0.001393"
$ws.Range("E15").Value = "-0.68%"
$ws.Range("G15").Value = "14"
$ws.Range("D16").Value = "0.006096"
$ws.Range("E16").Value = "-1.70%"
$ws.Range("G16").Value = "14"
$ws.Range("D17").Value = "3.729"
$ws.Range("E17").Value = "-3.17%"
$ws.Range("G17").Value = "14"
$ws.Range("E18").Value = "3.10%"
$ws.Range("G18").Value = "14"
$ws.Range("E19").Value = "1.61%"
$ws.Range("G19").Value = "14"
$ws.Range("D20").Value = "0.1312"
$ws.Range("E20").Value = "-2.09%"
$ws.Range("G20").Value = "14"
$ws.Range("D21").Value = "4.797"
$ws.Range("E21").Value = "0.01%"
$ws.Range("G21").Value = "14"
$ws.Range("E22").Value = "-3.69%"
$ws.Range("G22").Value = "14"
$ws.Range("D23").Value = "0.04355"
$ws.Range("E23").Value = "-1.51%"
$ws.Range("G23").Value = "14"
$ws.Range("D24").Value = "0.001231"
$ws.Range("E24").Value = "-0.17%"
$ws.Range("G24").Value = "14"
$ws.Range("D25").Value = "0.004119"
$ws.Range("E25").Value = "-15.18%"
$ws.Range("G25").Value = "14"
$ws.Range("G26").Value = "14"
$ws.Range("E27").Value = "0.12%"
$ws.Range("G27").Value = "14"
$ws.Range("G28").Value = "14"
$ws.Range("G29").Value = "14"
$ws.Range("G30").Value = "14"
$ws.Range("G31").Value = "14"
$ws.Range("G32").Value = "14"
$ws.Range("G33").Value = "14"
$ws.Range("G34").Value = "14"
$ws.Range("G35").Value = "14"
$ws.Range("G36").Value = "14"
$ws.Range("G37").Value = "14"
$ws.Range("G38").Value = "14"
$ws.Range("D39").Value = "0.02075"
$ws.Range("E39").Value = "7.89%"
$ws.Range("G39").Value = "14"
$ws.Range("D40").Value = "0.05028"
$ws.Range("E40").Value = "-2.27%"
$ws.Range("G40").Value = "14"
$ws.Range("D41").Value = "0.007489"
$ws.Range("E41").Value = "-0.74%"
$ws.Range("G41").Value = "14"
$ws.Range("D42").Value = "0.01012"
$ws.Range("E42").Value = "-0.05%"
$ws.Range("G42").Value = "14"
$ws.Range("D43").Value = "0.1347"
$ws.Range("E43").Value = "0.11%"
$ws.Range("G43").Value = "14"
$ws.Range("D44").Value = "0.002144"
$ws.Range("E44").Value = "1.55%"
$ws.Range("G44").Value = "14"
$ws.Range("D45").Value = "0.008839"
$ws.Range("E45").Value = "-10.37%"
$ws.Range("G45").Value = "14"
$ws.Range("D46").Value = "0.00006189"
$ws.Range("E46").Value = "-0.05%"
$ws.Range("G46").Value = "14"
$ws.Range("E47").Value = "0.06%"
$ws.Range("G47").Value = "14"
$ws.Range("D48").Value = "0.002823"
$ws.Range("G48").Value = "14"
$ws.Range("E49").Value = "28.00%"
$ws.Range("G49").Value = "14"
$ws.Range("D50").Value = "0.00002104"
$ws.Range("E50").Value = "0.06%"
$ws.Range("G50").Value = "14"
$ws.Range("E51").Value = "0.06%"
$ws.Range("G51").Value = "14"

$editRange.Style = "Normal"
